# Update the DESeq2-style results table (log2FoldChange / padj columns)
# with the refreshed values used for the re-generated volcano plot PNG.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C2").Value = -0.05124647226819731
$ws.Range("D2").Value = 0.1165340111754933
$ws.Range("C3").Value = 1.075908433645916
$ws.Range("D3").Value = 0.2562821896478944
$ws.Range("C4").Value = -0.862585156419838
$ws.Range("D4").Value = 0.6123848634545964
$ws.Range("C5").Value = -0.5915662200416927
$ws.Range("D5").Value = 0.5692498516767089
$ws.Range("C6").Value = 0.496115412931781
$ws.Range("D6").Value = 0.9776611205750612
$ws.Range("C7").Value = 0.8059997732670778
$ws.Range("D7").Value = 0.08373810287694738
$ws.Range("C8").Value = -0.960430868263573
$ws.Range("D8").Value = 0.05704740441803435
$ws.Range("C9").Value = -1.283276340877165
$ws.Range("D9").Value = 0.3729204079113255
$ws.Range("C10").Value = -1.910412512405695
$ws.Range("D10").Value = 0.5160798982943893
$ws.Range("C11").Value = -0.07993869337389428
$ws.Range("D11").Value = 0.5065168976828339
$ws.Range("C12").Value = -0.9222513109305691
$ws.Range("D12").Value = 0.2090149122343908
$ws.Range("C13").Value = -0.8010957157747601
$ws.Range("D13").Value = 0.6147907800365555
$ws.Range("C14").Value = 1.591884368780668
$ws.Range("D14").Value = 0.06501636056340265
$ws.Range("C15").Value = -0.3210274743165798
$ws.Range("D15").Value = 0.19669365233223
$ws.Range("C16").Value = 0.7173680415153602
$ws.Range("D16").Value = 0.06267052133251527
$ws.Range("C17").Value = 0.7245411693416838
$ws.Range("D17").Value = 0.1190185795011929
$ws.Range("C18").Value = -1.655342758868531
$ws.Range("D18").Value = 0.1228246884167578
$ws.Range("C19").Value = -0.3226804151994406
$ws.Range("D19").Value = 0.6615440383159179
$ws.Range("C20").Value = -0.2911080528624239
$ws.Range("D20").Value = 0.9841847084960202
$ws.Range("C21").Value = -1.530124231068745
$ws.Range("D21").Value = 0.6643174011264494
$ws.Range("C22").Value = 0.8804113352684855
$ws.Range("D22").Value = 0.5943418104231554
$ws.Range("C23").Value = -1.631203511293053
$ws.Range("D23").Value = 0.04702208276932096
$ws.Range("C24").Value = 0.1664752793604172
$ws.Range("D24").Value = 0.03769099897736761
$ws.Range("C25").ClearContents()
$ws.Range("D25").ClearContents()
$ws.Range("C26").Value = 0.1524179177153441
$ws.Range("D26").Value = 0.02689275956014825
$ws.Range("C27").Value = 1.845735829739362
$ws.Range("D27").Value = 0.1024339862552363
$ws.Range("C28").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("C29").Value = -0.3759048248361666
$ws.Range("D29").Value = 0.5243939033757079
$ws.Range("C30").Value = 1.611479909065487
$ws.Range("D30").Value = 0.6586129274711541
$ws.Range("C31").Value = -0.4397314157413
$ws.Range("D31").Value = 0.01672104550168652
$ws.Range("C32").Value = -1.402825564917863
$ws.Range("D32").Value = 0.7075175117503557
$ws.Range("C33").Value = 0.2210867955338335
$ws.Range("D33").Value = 0.9696923261399181
$ws.Range("C34").Value = -0.7420197696734002
$ws.Range("D34").Value = 0.4619784935688916
$ws.Range("C35").Value = -0.2750984841495288
$ws.Range("D35").Value = 0.6487584902811252
$ws.Range("C36").ClearContents()
$ws.Range("D36").ClearContents()
$ws.Range("C37").Value = 0.5603937532529546
$ws.Range("D37").Value = 0.8307877469683947
$ws.Range("C38").Value = -0.8242859282563622
$ws.Range("D38").Value = 0.9581095100560945
$ws.Range("C39").Value = 2.444120456206064
$ws.Range("D39").Value = 0.04638484900854045
$ws.Range("C40").Value = 0.7686132473480797
$ws.Range("D40").Value = 0.5562684433367024
$ws.Range("C41").Value = 1.374239120819454
$ws.Range("D41").Value = 0.5349678785378849
$ws.Range("C42").Value = 0.6350520075353919
$ws.Range("D42").Value = 0.8006116038739939
$ws.Range("C43").Value = -0.5915642819854271
$ws.Range("D43").Value = 0.7411920224712418
$ws.Range("C44").Value = 0.2704931499780999
$ws.Range("D44").Value = 0.8253233326976952
$ws.Range("C45").Value = -0.4462513212349941
$ws.Range("D45").Value = 0.2720473610571977
$ws.Range("C46").Value = -0.2381837668807433
$ws.Range("D46").Value = 0.01342513702925496
$ws.Range("C47").Value = -1.111909800864932
$ws.Range("D47").Value = 0.02363084715052655
$ws.Range("C48").Value = 0.9775695303819989
$ws.Range("D48").Value = 0.07164632373198765
$ws.Range("C49").Value = -0.3850449240827401
$ws.Range("D49").Value = 0.9190747315660529
$ws.Range("C50").Value = 1.02066401306132
$ws.Range("D50").Value = 0.550157945931504
$ws.Range("C51").Value = -0.4285865435911622
$ws.Range("D51").Value = 0.9896147455832123
$ws.Range("C52").Value = -0.3463420192076028
$ws.Range("D52").Value = 0.6433319544610455
$ws.Range("C53").Value = -0.9223454756643625
$ws.Range("D53").Value = 0.6445743966068982
$ws.Range("C54").ClearContents()
$ws.Range("D54").ClearContents()
$ws.Range("C55").Value = -0.1602341770834688
$ws.Range("D55").Value = 0.2985668078690291
$ws.Range("C56").Value = -0.6574329302513127
$ws.Range("D56").Value = 0.344312488365936
$ws.Range("C57").Value = -0.7706613668751215
$ws.Range("D57").Value = 0.8915941815112437
$ws.Range("C58").Value = 0.9145294978714179
$ws.Range("D58").Value = 0.34396456277867
$ws.Range("C59").Value = 0.1944381031280404
$ws.Range("D59").Value = 0.3699989987543775
$ws.Range("C60").Value = 1.430469462501794
$ws.Range("D60").Value = 0.6088044527644301
$ws.Range("C61").Value = -1.574250373019882
$ws.Range("D61").Value = 0.865040112261611
$ws.Range("C62").Value = -0.9636819929880518
$ws.Range("D62").Value = 0.3424253452754747
$ws.Range("C63").Value = -0.1053520452352895
$ws.Range("D63").Value = 0.9645360165064154
$ws.Range("C64").Value = 0.001385285828498836
$ws.Range("D64").Value = 0.6365474652221111
$ws.Range("C65").ClearContents()
$ws.Range("D65").ClearContents()
$ws.Range("C66").Value = 0.8289027458449272
$ws.Range("D66").Value = 0.8419375912535541
$ws.Range("C67").Value = -0.3566807741546517
$ws.Range("D67").Value = 0.2131188243232254
$ws.Range("C68").Value = -1.896545956671251
$ws.Range("D68").Value = 0.9746195414982775
$ws.Range("C69").Value = -0.05818895410316782
$ws.Range("D69").Value = 0.6423275673554977
$ws.Range("C70").Value = -1.754990556949109
$ws.Range("D70").Value = 0.002724282840635417
$ws.Range("C71").Value = 0.4033832801608584
$ws.Range("D71").Value = 0.5441148321379273
$ws.Range("C72").Value = 0.2088812391504386
$ws.Range("D72").Value = 0.5753301617450718
$ws.Range("C73").Value = 1.4396483170419
$ws.Range("D73").Value = 0.3138981979377987
$ws.Range("C74").Value = -1.599052686908847
$ws.Range("D74").Value = 0.1968248684760824
$ws.Range("C75").ClearContents()
$ws.Range("D75").ClearContents()
$ws.Range("C76").Value = 1.556366527463733
$ws.Range("D76").Value = 0.5858971137754244
$ws.Range("C77").Value = 0.09205738119948087
$ws.Range("D77").Value = 0.03099507291631654
$ws.Range("C78").Value = -0.1405579790414384
$ws.Range("D78").Value = 0.0709138775051058
$ws.Range("C79").Value = -0.541458759951653
$ws.Range("D79").Value = 0.5588197190438194
$ws.Range("C80").Value = -0.4386298772198119
$ws.Range("D80").Value = 0.339486897531738
$ws.Range("C81").ClearContents()
$ws.Range("D81").ClearContents()
$ws.Range("C82").Value = -1.770940778519035
$ws.Range("D82").Value = 0.2627873695400764
$ws.Range("C83").Value = -0.5930212646339734
$ws.Range("D83").Value = 0.8136768373783398
$ws.Range("C84").Value = -1.359318862453168
$ws.Range("D84").Value = 0.8206599874585043
$ws.Range("C85").Value = 0.5542485346813597
$ws.Range("D85").Value = 0.3550442399372543
$ws.Range("C86").Value = -1.329515986631947
$ws.Range("D86").Value = 0.5997861353974804
$ws.Range("C87").Value = 0.3181974179141019
$ws.Range("D87").Value = 0.2987334853664169
$ws.Range("C88").Value = -2.514551893717444
$ws.Range("D88").Value = 0.8960125092786448
$ws.Range("C89").Value = -2.137024664737153
$ws.Range("D89").Value = 0.2618206529954723
$ws.Range("C90").Value = 1.769882961555094
$ws.Range("D90").Value = 0.2279298648862725
$ws.Range("C91").Value = 2.3629718871612
$ws.Range("D91").Value = 0.3658031577842265
$ws.Range("C92").Value = 0.1599691707801534
$ws.Range("D92").Value = 0.6220653903370444
$ws.Range("C93").Value = -0.7485677948369388
$ws.Range("D93").Value = 0.8852216133753907
$ws.Range("C94").ClearContents()
$ws.Range("D94").ClearContents()
$ws.Range("C95").Value = -1.398782842501166
$ws.Range("D95").Value = 0.9878382817221761
$ws.Range("C96").Value = 1.226975111151378
$ws.Range("D96").Value = 0.5022998301901366
$ws.Range("C97").ClearContents()
$ws.Range("D97").ClearContents()
$ws.Range("C98").Value = 0.172941246157735
$ws.Range("D98").Value = 0.1673830032250285
$ws.Range("C99").ClearContents()
$ws.Range("D99").ClearContents()
$ws.Range("C100").Value = 1.415584947771703
$ws.Range("D100").Value = 0.9726268232037465
$ws.Range("C101").Value = -1.479720995396407
$ws.Range("D101").Value = 0.3922692406903486
